$p = $ppt.ActivePresentation

# --- Slide 2 (Security office area): rename "PROPERTY OFFICE" -> "STORAGE ROOM" ---
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt.Trim() -eq "PROPERTY OFFICE") {
            $shp.TextFrame.TextRange.Text = "STORAGE ROOM"
        }
    }
}

# --- Slide 4: add a new "NO CAMERA FEED" room-name textbox, matching the style ---
# of the existing "SHOWER BLOCK" / "TOILET BLOCK" labels already on that slide.
# Duplicating an existing label keeps all the text/shape formatting (font size,
# shadow effect, no-fill box, centered/no-wrap body) identical, as in the diff.
$s4 = $p.Slides.Item(4)
$template = $null
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $shp = $s4.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text.Trim() -eq "SHOWER BLOCK") {
        $template = $shp
    }
}
if ($template -eq $null) {
    $template = $s4.Shapes.Item($s4.Shapes.Count)
}

$dup = $template.Duplicate()
$newShape = $dup.Item(1)
$newShape.Name = "TextBox 1"
$newShape.Left = 0
$newShape.Top = 226.77173228346456
$newShape.TextFrame.TextRange.Text = "NO CAMERA FEED"
